# Project "Sample Project" is saved: the rule-table cell B11 (row "R40"),
# which used to hold the label "R40", now holds the text "1".
#
# A leading apostrophe is used so Excel stores the numeric-looking text
# "1" as a literal string (shared-string / text cell) rather than
# re-interpreting it as the number 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
